$d = $word.ActiveDocument

$lq = [char]8220
$rq = [char]8221
$apos = [char]8217

# ---------------------------------------------------------------------------
# 1. Intro paragraph: "turing" -> "Turing" and rework the Busy Beaver sentence.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This project consisted of two portions: building a turing machine simulator and creating a so-called " + $lq + "Busy Beaver" + $rq + " machine that writes the most non-zeros to an infinite 2-way tape given blank input before halting.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This project consisted of two portions: building a Turing machine simulator and creating a so-called " + $lq + "Busy Beaver" + $rq + " machine that maximizes the sum of non-zeros to an infinite 2-way tape, given blank input before halting.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2. "First, we will discuss the design of our turing machine simulator"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "First, we will discuss the design of our turing machine simulator",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "First, we will discuss the design of our Turing machine simulator",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 3. "1. Checks if the current is a the halting state" -> remove stray "a "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "1. Checks if the current is a the halting state",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Checks if the current is the halting state",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 4. Typo fix: "Se the next transition" -> "Set the next transition"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Se the next transition for the new state and repeat the algorithm.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Set the next transition for the new state and repeat the algorithm.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 5. Expand the testing paragraph into three paragraphs.
# ---------------------------------------------------------------------------
$oldTesting = "First, an input file matching the one given on Blackboard was created and stepped through to ensure the simulator worked properly."
$newA = "First, we tested our parser by parsing the example input files and making sure the resulting Machine object actually matched the structure described in the file. "
$newB = "Then, an input file matching the one given on Blackboard was created and stepped through to ensure the simulator worked properly and actually halted. "
$newC = "We then manually simulated a few simple Turing machines, and simulated them with our program, to make sure the output was correct."
$replacement = $newA + "^p" + $newB + "^p" + $newC

$d.Content.Find.Execute($oldTesting, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null

# Locate paragraph C ("We then manually simulated...") so we can anchor the
# bookmark and the new appended paragraphs after it.
$idxC = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq $newC) {
        $idxC = $i
    }
}

# Word leaves a "_GoBack" bookmark at the last edit location.
$d.Bookmarks.Add("_GoBack", $d.Paragraphs($idxC).Range) | Out-Null

# ---------------------------------------------------------------------------
# 6. Append the new "Optimizing the TM simulator" section.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

$pHeading = $d.Paragraphs.Add($lastPara.Range)
$pHeading.Range.Text = "Optimizing the TM simulator"
$pHeading.Range.Font.Bold = 1
$pHeading.Range.Font.Italic = 1

$pE = $d.Paragraphs.Add($lastPara.Range)
$pE.Range.Text = "Once we had a functioning Turing machine simulator, we began work on optimizing it, since it was pretty slow (our best Busy Beaver machine took almost 20 minutes to run). To bring this run time down to a more reasonable level, we looked at where the speed bottleneck was happening, which was outputting the contents of the tape. We realized that this output was actually happening while the Turing machine was being simulated, rather than after it finished, which was not only slow, but also was likely giving us extra and incorrect output. By simply waiting to output the tape contents until after the simulation finishes, we were able to reduce the run time from 20 minutes to about 6 seconds (!!!)."

$pF = $d.Paragraphs.Add($lastPara.Range)
$pF.Range.Text = "Looking through the code, we also saw a few opportunities to make the simulator more efficient. The biggest inefficiency we noticed was that we were using an ArrayList<Integer> for our tape. Since we simulated the tape by inserting new tape cells at the beginning and end of the ArrayList, the whole tape needed shifted in order to move the tape head left, in some cases, and occasional ArrayList resizes meant lots of data was moving around unnecessarily. While most tape operations were O(1), every once in a while, they" + $apos + "d be O(n). To speed this up, we implemented the tape as a simple linked list of ints (simple, because there" + $apos + "s never any need to insert new tape cells anywhere except the beginning and end of the linked list, and deletions never happen). By switching to this implementation, we were able to reduce every tape operation, except printing the tape contents, to O(1). "

$pG = $d.Paragraphs.Add($lastPara.Range)
$pG.Range.Text = "Also, but using primitive ints, rather than Integers, we were able to avoid the potential overhead from the JVM constantly wrapping/unwrapping Integer objects. Similarly, by switching from using primitive chars, rather than String objects, from some data, we found small gains in performance."

$pH = $d.Paragraphs.Add($lastPara.Range)
$pH.Range.Text = "With these optimizations, we were able to reduce our 6 second run time by nearly half to about 3.5 seconds."

# ---------------------------------------------------------------------------
# 7. Remove the old trailing empty paragraph (the new content now ends the
#    document body immediately before the sectPr).
# ---------------------------------------------------------------------------
$delRange = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End)
$delRange.Delete() | Out-Null

Write-Host "done"
